# "updated legacy GSC export data"
#
# The legacy GSC export table ("Table" sheet, first / active sheet in the
# workbook) contained a stray leading row for 2025-11-02 whose
# "No video indexed" / "Video indexed" values were blank (the export run
# for that day never populated them). That row is removed, which shifts
# the remaining 88 daily rows (2025-11-03 .. 2026-01-29) up by one.
#
# Additionally the last row in the table (2026-01-29) had its
# "Impressions" cell stored as an empty string instead of the numeric 0
# used by every other row - that gets corrected to a real number 0 once
# it becomes the new final row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete 2025-11-02 row (row 2, right under the header row).
$ws.Rows.Item(2).Delete()

# After the delete, the table runs from row 2 to this last row.
$lastRow = $ws.UsedRange.Rows.Count

# Fix the trailing Impressions value (column D) for the now-last row
# (2026-01-29) so it is a proper number instead of an empty string.
$ws.Range("D" + $lastRow).Value = 0
